$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(2, 2).Range.Text = "422"
$tbl.Cell(3, 2).Range.Text = "388 (95.1)"
$tbl.Cell(5, 2).Range.Text = "380 (93.1)"
$tbl.Cell(6, 2).Range.Text = "86 (20.9)"
$tbl.Cell(7, 2).Range.Text = "193 (51.7)"
